$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Titus" (row group for rows 18-22) evaluated the seq2seq ("human_s2s_lstm_att_long_win-*")
# outputs whose filenames already occupy header columns AA2:AF2. Fill in his
# scores (Fluency, Coherence, Meaningfulness, Poeticness, Overall) for the six
# seq2seq files - all scored 1/5.
$cols = @("AA", "AB", "AC", "AD", "AE", "AF")
for ($r = 18; $r -le 22; $r++) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = 1
    }
}
